# Apply text replacements based on the diff
$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-05-25 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-26 Sunday", 2) | Out-Null
$d.Content.Find.Execute("35×35=1225", $true, $false, $false, $false, $false, $true, 1, $false, "50×12=600", 2) | Out-Null
$d.Content.Find.Execute("51×82=4182", $true, $false, $false, $false, $false, $true, 1, $false, "19×39=741", 2) | Out-Null
$d.Content.Find.Execute("18×82=1476", $true, $false, $false, $false, $false, $true, 1, $false, "51×58=2958", 2) | Out-Null
$d.Content.Find.Execute("37×41=1517", $true, $false, $false, $false, $false, $true, 1, $false, "35×37=1295", 2) | Out-Null
$d.Content.Find.Execute("84×71=5964", $true, $false, $false, $false, $false, $true, 1, $false, "17×66=1122", 2) | Out-Null
$d.Content.Find.Execute("91×73=6643", $true, $false, $false, $false, $false, $true, 1, $false, "46×53=2438", 2) | Out-Null
$d.Content.Find.Execute("22×65=1430", $true, $false, $false, $false, $false, $true, 1, $false, "87×60=5220", 2) | Out-Null
$d.Content.Find.Execute("45×70=3150", $true, $false, $false, $false, $false, $true, 1, $false, "12×54=648", 2) | Out-Null
$d.Content.Find.Execute("47×92=4324", $true, $false, $false, $false, $false, $true, 1, $false, "39×32=1248", 2) | Out-Null
$d.Content.Find.Execute("94×83=7802", $true, $false, $false, $false, $false, $true, 1, $false, "80×37=2960", 2) | Out-Null
$d.Content.Find.Execute("98×55=5390", $true, $false, $false, $false, $false, $true, 1, $false, "72×32=2304", 2) | Out-Null
$d.Content.Find.Execute("30×23=690", $true, $false, $false, $false, $false, $true, 1, $false, "91×20=1820", 2) | Out-Null
$d.Content.Find.Execute("46×59=2714", $true, $false, $false, $false, $false, $true, 1, $false, "64×85=5440", 2) | Out-Null
$d.Content.Find.Execute("20×89=1780", $true, $false, $false, $false, $false, $true, 1, $false, "43×36=1548", 2) | Out-Null
$d.Content.Find.Execute("70×11=770", $true, $false, $false, $false, $false, $true, 1, $false, "59×89=5251", 2) | Out-Null
$d.Content.Find.Execute("37×93=3441", $true, $false, $false, $false, $false, $true, 1, $false, "90×35=3150", 2) | Out-Null
$d.Content.Find.Execute("27×25=675", $true, $false, $false, $false, $false, $true, 1, $false, "48×16=768", 2) | Out-Null
$d.Content.Find.Execute("97×45=4365", $true, $false, $false, $false, $false, $true, 1, $false, "66×30=1980", 2) | Out-Null
$d.Content.Find.Execute("14×42=588", $true, $false, $false, $false, $false, $true, 1, $false, "74×34=2516", 2) | Out-Null
$d.Content.Find.Execute("97×26=2522", $true, $false, $false, $false, $false, $true, 1, $false, "31×31=961", 2) | Out-Null
$d.Content.Find.Execute("28×29=812", $true, $false, $false, $false, $false, $true, 1, $false, "84×77=6468", 2) | Out-Null
$d.Content.Find.Execute("16×40=640", $true, $false, $false, $false, $false, $true, 1, $false, "22×99=2178", 2) | Out-Null
$d.Content.Find.Execute("80×46=3680", $true, $false, $false, $false, $false, $true, 1, $false, "92×21=1932", 2) | Out-Null
$d.Content.Find.Execute("34×58=1972", $true, $false, $false, $false, $false, $true, 1, $false, "78×93=7254", 2) | Out-Null
$d.Content.Find.Execute("83×44=3652", $true, $false, $false, $false, $false, $true, 1, $false, "27×68=1836", 2) | Out-Null
